# Adds ML Engineer (Recommendation Systems) and ML Engineer (Robot Learning) job rows
# to the 'Job Listings' sheet, rows 70-86, columns A-I.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('TikTok', 'Machine Learning Engineer, TikTok Recommendation', 'Big Tech', 'San Jose, CA', 'ML Engineer (Recommendation Systems)', 'https://careers.tiktok.com/position/7112323775379343623/detail', 'Master''s or PhD in Computer Science, Machine Learning, or related field; experience with large-scale recommendation systems; strong publication record preferred', 'Build large-scale (10M-100M) recommendation algorithms and systems; design user interest models for short and long-term preferences; develop CTR/CVR prediction models; build real-time data pipelines and feature engineering', 'Python, C++, PyTorch/TensorFlow, distributed systems, deep learning, recommendation algorithms, real-time ML systems, feature engineering, A/B testing')
    ,@('Amazon', 'Lead AI/ML Engineer, Recommendation Systems (Prime Video)', 'Big Tech', 'Seattle, WA', 'ML Engineer (Recommendation Systems)', 'https://amazon.jobs/en/jobs/3065041/lead-ai-ml-engineer-recommendation-systems', 'PhD or Master''s degree in CS, ML, or related field; 6+ years applied research experience; 5+ years building ML models for business applications', 'Develop AI solutions for personalization; design and conduct offline and online A/B experiments; stay current with latest modeling techniques; lead technical direction for recommendation systems', 'Python, PyTorch/TensorFlow, deep learning, recommendation systems, A/B testing, causal inference, distributed computing, AWS, Spark')
    ,@('Amazon', 'Machine Learning Engineer - Personalization', 'Big Tech', 'Seattle, WA', 'ML Engineer (Recommendation Systems)', 'https://www.amazon.jobs/en/jobs/3065042/machine-learning-engineer-personalization', 'Master''s or PhD in CS, CE, ML or related field; 3+ years building ML models for business applications; experience with recommendation and search systems', 'Develop ML models for recommendation and search systems using deep learning, online learning, and optimization; tackle cold-start problems with foundation models; build scalable personalization systems', 'Python, PyTorch/TensorFlow, deep learning, online learning, optimization, recommendation systems, foundation models, AWS, distributed systems')
    ,@('Amazon', 'Machine Learning Engineer, Amazon Music ML & Personalization', 'Big Tech', 'San Francisco, CA', 'ML Engineer (Recommendation Systems)', 'https://www.amazon.jobs/en/jobs/3004439/machine-learning-engineer-amazon-music-ml-personalization', 'BS/MS in Computer Science, ML, or related field; experience with recommendation systems and personalization; familiarity with music/audio domain a plus', 'Deliver personalized recommendations for millions of Amazon Music customers globally; own models underlying all recommendation surfaces across mobile, web, and Alexa; build data pipelines and approximate nearest neighbor services', 'Python, PyTorch/TensorFlow, recommendation systems, approximate nearest neighbors, data pipelines, A/B testing, AWS, distributed systems, audio/music ML')
    ,@('Meta', 'Senior Machine Learning Engineer, Recommendation Systems', 'FAANG', 'Menlo Park, CA', 'ML Engineer (Recommendation Systems)', 'https://www.metacareers.com/jobs/2472910023050061', 'BS/MS/PhD in Computer Science, ML, or related field; 5+ years experience with recommendation systems; experience at companies like Google, Meta, or Twitter preferred', 'Architect, design, implement, and test large-scale distributed recommendation systems end-to-end; lead team building recommendation product for Short Videos in Facebook Newsfeed; own inventory design, content retrieval, ranking stages, delivery, and ML infra', 'Python, C++, PHP/Hack, PyTorch, TensorFlow, Spark, Presto, distributed systems, ranking systems, content retrieval, real-time ML')
    ,@('TikTok', 'ML Engineer Graduate - Multimodal Recommendation (2026 Start)', 'Big Tech', 'San Jose, CA', 'ML Engineer (Recommendation Systems)', 'https://jobright.ai/jobs/info/68cc15b916d00d2beeb19757', 'BS/MS in Computer Science, ML, or related field; graduating 2026; research experience in multimodal learning or recommendation systems preferred', 'Lead multimodal algorithm development for short-video business; explore applications of multimodal technologies in recommendation systems; conduct cutting-edge research in multimodal and MLLM technologies', 'Python, PyTorch/TensorFlow, multimodal learning, LLMs, MLLMs, recommendation systems, computer vision, NLP, deep learning, research')
    ,@('Netflix', 'Machine Learning Engineer - Personalization', 'Big Tech', 'Los Gatos, CA', 'ML Engineer (Recommendation Systems)', 'https://explore.jobs.netflix.net/careers/job/790312415414', 'MS/PhD in Computer Science, ML, Statistics, or related field; experience with recommendation systems and personalization at scale', 'Design, build, and deploy ML models powering personalized recommendations; develop scalable ML pipelines and infrastructure; work on content optimization, search, and streaming quality', 'Python, PyTorch/TensorFlow, recommendation systems, personalization, distributed systems, A/B testing, causal inference, data pipelines, Spark')
    ,@('Bedrock Robotics', 'Machine Learning Engineer: Robotics ML Generalist', 'Startup', 'San Francisco, CA / New York, NY', 'ML Engineer (Robot Learning)', 'https://see.stanford.edu/course/cs229', 'Hands-on robotics deployment experience; deep learning expertise for perception and behavior systems', 'Deploy ML models to robots; build data labeling pipelines; real-time optimization for robotic systems', 'PyTorch, Python, C++/Rust, ONNX/TensorRT, sensor data processing (camera/lidar/IMU)')
    ,@('Scale AI', 'Machine Learning Research Engineer - Robotics', 'Big Tech', 'San Francisco, CA', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', '3+ years ML experience; ability to take research to production; mid/senior level; H1B possible', 'Research and develop ML for robotics applications; bridge research to production systems', 'PyTorch/TensorFlow, reinforcement learning, robotics simulation, research methods')
    ,@('Toyota Research Institute', 'Robotics ML Engineer - Platforms for Vision Language Action', 'Big Tech', 'Los Altos, CA', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', 'Mid/senior level; $176K-$264K salary range; experience with foundation models and robotics platforms', 'Build and evaluate learning methods for manipulation; work with Vision-Language-Action (VLA) models', 'Foundation models, VLA models, robotics platforms, manipulation learning, PyTorch')
    ,@('Johnson & Johnson MedTech', 'Staff Machine Learning Engineer, Robotics', 'Big Tech', 'Santa Clara, CA', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', '6+ years experience; senior/staff level; $141K-$228K; hybrid work; H1B sponsorship', 'Deploy ML models for medical robotics; computer vision for surgical applications; edge computing optimization', 'ML deployment, computer vision, edge computing, medical robotics, real-time systems')
    ,@('Berkshire Grey', 'Principal Robotics Engineer, Machine Learning', 'Big Tech', 'Bedford, MA', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', '6+ years robotics/ML experience; senior/lead/staff level; $162K-$200K; near Boston', 'Build ML systems for warehouse robotics; scale production autonomy systems', 'ML for robotics, autonomy systems, production scaling, warehouse automation')
    ,@('Pickle Robot Company', 'Perception and Computer Vision Engineer', 'Startup', 'Boston, MA', 'ML Engineer (Robot Learning)', 'https://www.linkedin.com/jobs/robotics-engineer-jobs-boston-ma', 'Local to Boston area; experience with robotics sensor integration', 'Develop computer vision and ML for perception; real-time processing for robotic systems', 'Computer vision, ML for perception, robotics sensors, real-time processing, Python')
    ,@('Diligent Robotics', 'Sr ML Engineer, Robotics (Perception)', 'Startup', 'Austin, TX', 'ML Engineer (Robot Learning)', 'https://www.linkedin.com/jobs/deep-learning-robotics-jobs', '5+ years ML/CV experience; senior level; edge deployment expertise', 'Build perception systems for healthcare robots; detection and segmentation models; edge deployment', 'PyTorch/TensorFlow, detection/segmentation, edge deployment, healthcare robotics')
    ,@('XPENG', 'Machine Learning Engineer - LLM, AI & Robotics', 'Big Tech', 'Santa Clara, CA', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', 'Senior level; $149K-$252K; automotive robotics focus', 'Develop LLM and AI systems for automotive robotics applications', 'LLMs, robotics AI, deep learning, automotive systems, autonomous vehicles')
    ,@('Path Robotics', 'Machine Learning Engineer, Robot Learning', 'Startup', 'Columbus, OH', 'ML Engineer (Robot Learning)', 'https://jobright.ai/jobs/machine-learning-engineer---robotics-jobs-in-united-states', 'Mid-level; H1B sponsorship available; industrial automation focus', 'Build ML systems for industrial robot automation; imitation and reinforcement learning', 'Robot learning, imitation learning, reinforcement learning, industrial automation')
    ,@('RobCo', 'Machine Learning Engineer - Robot Learning', 'Startup', 'Europe/US (Remote possible)', 'ML Engineer (Robot Learning)', 'https://robco.jobs.personio.de/job/2450134?language=en', 'ML/robotics intersection experience; ability to deploy to hardware', 'Develop learning-based manipulation methods; deploy to real robot hardware; cross-team collaboration', 'Learning-based manipulation, systems engineering, hardware deployment, sim-to-real transfer')
)

$startRow = 70
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}

Write-Output "Added $($data.Count) rows starting at row $startRow"
